$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel auto-converting
# numeric-looking strings (e.g. "1.001") into real numbers. We build a formula
# that evaluates to the literal text, then Copy + PasteSpecial(values-only) to
# flatten it back down to a plain inline string cell (no leftover formula,
# no style change).
function Set-TextValue($ws, $row, $col, $text) {
    $escaped = $text.Replace('"', '""')
    $cell = $ws.Cells.Item($row, $col)
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

Set-TextValue $ws 2 4 '30.025.02'
$ws.Cells.Item(2, 5).Value = '  -0.85%  '

Set-TextValue $ws 3 4 '1.906.13'
$ws.Cells.Item(3, 5).Value = '  -0.75%  '

Set-TextValue $ws 4 4 '1.001'
$ws.Cells.Item(4, 5).Value = '  -0.26%  '

Set-TextValue $ws 5 4 '0.7577'
$ws.Cells.Item(5, 5).Value = '  +2.21%  '

Set-TextValue $ws 6 4 '241.27'
$ws.Cells.Item(6, 5).Value = '  -1.21%  '

Set-TextValue $ws 7 4 '1.000'
$ws.Cells.Item(7, 5).Value = '  -0.27%  '

Set-TextValue $ws 8 4 '0.3078'
$ws.Cells.Item(8, 5).Value = '  -1.94%  '

Set-TextValue $ws 9 4 '25.48'
$ws.Cells.Item(9, 5).Value = '  -6.57%  '

Set-TextValue $ws 10 4 '0.06892'
$ws.Cells.Item(10, 5).Value = '  -1.15%  '

Set-TextValue $ws 11 4 '0.08004'
$ws.Cells.Item(11, 5).Value = '  +0.16%  '

Set-TextValue $ws 12 4 '0.7540'
$ws.Cells.Item(12, 5).Value = '  -2.46%  '

Set-TextValue $ws 13 4 '1.902.65'
$ws.Cells.Item(13, 5).Value = '  -1.24%  '

Set-TextValue $ws 14 4 '5.246'
$ws.Cells.Item(14, 5).Value = '  -1.15%  '

Set-TextValue $ws 15 4 '91.64'
$ws.Cells.Item(15, 5).Value = '  -0.06%  '

Set-TextValue $ws 16 4 '6.186'
$ws.Cells.Item(16, 5).Value = '  +5.76%  '

Set-TextValue $ws 17 4 '30.024.52'
$ws.Cells.Item(17, 5).Value = '  -0.97%  '

Set-TextValue $ws 18 4 '14.01'
$ws.Cells.Item(18, 5).Value = '  -1.65%  '

Set-TextValue $ws 19 4 '0.000007715'
$ws.Cells.Item(19, 5).Value = '  -1.81%  '

Set-TextValue $ws 20 4 '237.04'
$ws.Cells.Item(20, 5).Value = '  -3.69%  '

Set-TextValue $ws 21 4 '1.000'
$ws.Cells.Item(21, 5).Value = '  -0.25%  '

Set-TextValue $ws 22 4 '2.157.49'
$ws.Cells.Item(22, 5).Value = '  -0.89%  '

Set-TextValue $ws 23 4 '1.001'
$ws.Cells.Item(23, 5).Value = '  -0.22%  '

Set-TextValue $ws 24 4 '7.025'
$ws.Cells.Item(24, 5).Value = '  +5.33%  '

Set-TextValue $ws 25 4 '9.295'
$ws.Cells.Item(25, 5).Value = '  -1.34%  '

Set-TextValue $ws 26 4 '166.16'
$ws.Cells.Item(26, 5).Value = '  +0.57%  '

Set-TextValue $ws 27 4 '18.82'
$ws.Cells.Item(27, 5).Value = '  -0.69%  '

Set-TextValue $ws 28 4 '0.1299'
$ws.Cells.Item(28, 5).Value = '  +2.37%  '

Set-TextValue $ws 29 4 '2.067'
$ws.Cells.Item(29, 5).Value = '  -3.22%  '

$ws.Cells.Item(30, 5).Value = '  -1.14%  '

Set-TextValue $ws 31 4 '1.522'
$ws.Cells.Item(31, 5).Value = '  -1.83%  '

Set-TextValue $ws 32 4 '4.305'
$ws.Cells.Item(32, 5).Value = '  -1.03%  '

Set-TextValue $ws 33 4 '4.043'
$ws.Cells.Item(33, 5).Value = '  -0.79%  '

Set-TextValue $ws 34 4 '0.05462'
$ws.Cells.Item(34, 5).Value = '  +5.36%  '

Set-TextValue $ws 35 4 '1.288'
$ws.Cells.Item(35, 5).Value = '  -1.06%  '

Set-TextValue $ws 36 4 '0.7366'
$ws.Cells.Item(36, 5).Value = '  -1.76%  '

Set-TextValue $ws 37 4 '2.715'
$ws.Cells.Item(37, 5).Value = '  -2.12%  '

Set-TextValue $ws 38 4 '0.01939'
$ws.Cells.Item(38, 5).Value = '  -0.16%  '

$ws.Cells.Item(39, 5).Value = '  -0.83%  '

Set-TextValue $ws 40 4 '6.241'
$ws.Cells.Item(40, 5).Value = '  -2.58%  '

Set-TextValue $ws 41 4 '0.4446'
$ws.Cells.Item(41, 5).Value = '  -0.32%  '

Set-TextValue $ws 42 4 '72.77'
$ws.Cells.Item(42, 5).Value = '  -4.22%  '

Set-TextValue $ws 43 4 '1.942'
$ws.Cells.Item(43, 5).Value = '  -0.27%  '

Set-TextValue $ws 44 4 '1.000'
$ws.Cells.Item(44, 5).Value = '  -0.16%  '

Set-TextValue $ws 45 4 '0.8307'
$ws.Cells.Item(45, 5).Value = '  -0.92%  '

Set-TextValue $ws 46 4 '7.657'
$ws.Cells.Item(46, 5).Value = '  -0.20%  '

Set-TextValue $ws 47 4 '101.49'
$ws.Cells.Item(47, 5).Value = '  +0.16%  '

Set-TextValue $ws 48 4 '9.847'
$ws.Cells.Item(48, 5).Value = '  -0.14%  '

Set-TextValue $ws 49 4 '2.058.59'
$ws.Cells.Item(49, 5).Value = '  -1.70%  '

Set-TextValue $ws 50 4 '36.50'
$ws.Cells.Item(50, 5).Value = '  -1.68%  '

$ws.Cells.Item(51, 5).Value = '  -4.83%  '

$wb.Save()
